# 🚌 141: 30/12 23:35 LP1912+6203+6173
# Appends the latest scrape batch (run at 20:35:43/48/54) to the three
# tracking sheets and refreshes the "Última actualización" / "Total filas"
# header cells on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "LP1912"  (columns: A=meta, B=Hora_Scrap, C=Hora_Llegada,
#                      D=Línea, E=Minutos, F=Parada, G=Fecha)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 20:35:54"
$ws1.Cells.Item(3, 1).Value = "Total filas: 630"

$rows1 = @(
    @(617, "20:35:43", "20:46", "16_SANTA ANA",                11),
    @(618, "20:35:43", "20:52", "15_ABASTO",                   17),
    @(619, "20:35:43", "20:57", "23_HERNANDEZ",                22),
    @(620, "20:35:43", "21:01", "23_HERNANDEZ",                26),
    @(621, "20:35:43", "21:04", "84_COLONIA URQUIZA-ESC 49",   29),
    @(622, "20:35:43", "21:07", "215B_EL PATO",                32),
    @(623, "20:35:43", "21:17", "16_SANTA ANA",                42),
    @(624, "20:35:43", "21:20", "26_HERNANDEZ",                45),
    @(625, "20:35:43", "21:22", "15_ABASTO",                   47),
    @(626, "20:35:43", "21:31", "16_SANTA ANA",                56),
    @(627, "20:35:43", "21:32", "23_HERNANDEZ",                57),
    @(628, "20:35:43", "21:37", "17_ROMERO",                   62),
    @(629, "20:35:43", "21:47", "215A_EL PATO",                72),
    @(630, "20:35:43", "21:51", "10_OLMOS",                    76),
    @(631, "20:35:43", "22:07", "17_ROMERO",                   92)
)

foreach ($r in $rows1) {
    $rowNum = $r[0]
    $ws1.Cells.Item($rowNum, 2).Value = $r[1]
    $ws1.Cells.Item($rowNum, 3).Value = $r[2]
    $ws1.Cells.Item($rowNum, 4).Value = $r[3]
    $ws1.Cells.Item($rowNum, 5).Value = $r[4]
    $ws1.Cells.Item($rowNum, 6).Value = "LP1912"
    $ws1.Cells.Item($rowNum, 7).Value = "30/12/2025"
}

# ---------------------------------------------------------------------------
# Sheet 2: "LP1912-215"  (columns: A=meta, B=Fecha, C=Hora_Scrap,
#                          D=Hora_Llegada, E=Línea, F=Minutos, G=Parada)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 20:35:54"
$ws2.Cells.Item(3, 1).Value = "Total filas: 46"

$rows2 = @(
    @(46, "20:35:43", "21:07", "215B_EL PATO", 32),
    @(47, "20:35:43", "21:47", "215A_EL PATO", 72)
)

foreach ($r in $rows2) {
    $rowNum = $r[0]
    $ws2.Cells.Item($rowNum, 2).Value = "30/12/2025"
    $ws2.Cells.Item($rowNum, 3).Value = $r[1]
    $ws2.Cells.Item($rowNum, 4).Value = $r[2]
    $ws2.Cells.Item($rowNum, 5).Value = $r[3]
    $ws2.Cells.Item($rowNum, 6).Value = $r[4]
    $ws2.Cells.Item($rowNum, 7).Value = "LP1912"
}

# ---------------------------------------------------------------------------
# Sheet 3: "6203-6173"  (columns: A=meta, B=Fecha, C=Hora_Scrap,
#                         D=Hora_Llegada, E=Línea, F=Minutos, G=Parada)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 30/12/2025 20:35:54"
$ws3.Cells.Item(3, 1).Value = "Total filas: 76"

$ws3.Cells.Item(76, 2).Value = "30/12/2025"
$ws3.Cells.Item(76, 3).Value = "20:35:48"
$ws3.Cells.Item(76, 4).Value = "21:28"
$ws3.Cells.Item(76, 5).Value = "215C_LA PLATA"
$ws3.Cells.Item(76, 6).Value = 53
$ws3.Cells.Item(76, 7).Value = "L6203"

$ws3.Cells.Item(77, 2).Value = "30/12/2025"
$ws3.Cells.Item(77, 3).Value = "20:35:54"
$ws3.Cells.Item(77, 4).Value = "22:04"
$ws3.Cells.Item(77, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(77, 6).Value = 89
$ws3.Cells.Item(77, 7).Value = "L6173"
